$wb = $excel.ActiveWorkbook

# --- Sheet: strategy_id-0 (sheet1) ---
$ws = $wb.Worksheets.Item("strategy_id-0")

# Row 10: Y10:AS10 from 1 -> 2700000000
$ws.Range("Y10:AS10").Value = 2700000000

# Row 102: J102:AS102 -> 0
$ws.Range("J102:AS102").Value = 0

# Row 137: J137:AS137 -> 0.001
$ws.Range("J137:AS137").Value = 0.001

# Column C adjustments (decrement by 1) for sheet1
$ws.Range("C24").Value = 15
$ws.Range("C26").Value = 15
$ws.Range("C27").Value = 15
$ws.Range("C28").Value = 13
$ws.Range("C30").Value = 13
$ws.Range("C31").Value = 13
$ws.Range("C34").Value = 15
$ws.Range("C35").Value = 13
$ws.Range("C36").Value = 13
$ws.Range("C44").Value = 15
$ws.Range("C49").Value = 15
$ws.Range("C54").Value = 15
$ws.Range("C55").Value = 15
$ws.Range("C57").Value = 15
$ws.Range("C58").Value = 15
$ws.Range("C59").Value = 13
$ws.Range("C61").Value = 13
$ws.Range("C72").Value = 13
$ws.Range("C73").Value = 13
$ws.Range("C76").Value = 13
$ws.Range("C77").Value = 13
$ws.Range("C78").Value = 13
$ws.Range("C79").Value = 13
$ws.Range("C88").Value = 14
$ws.Range("C89").Value = 14
$ws.Range("C90").Value = 14
$ws.Range("C114").Value = 12
$ws.Range("C115").Value = 25
$ws.Range("C116").Value = 25
$ws.Range("C119").Value = 25
$ws.Range("C120").Value = 25

# --- Sheet: strategy_id-6004 ---
$ws = $wb.Worksheets.Item("strategy_id-6004")
$ws.Range("C2").Value = 15
$ws.Range("C3").Value = 15
$ws.Range("C4").Value = 15
$ws.Range("C5").Value = 13
$ws.Range("C6").Value = 13
$ws.Range("C7").Value = 13
$ws.Range("C8").Value = 15
$ws.Range("C9").Value = 13
$ws.Range("C10").Value = 13
$ws.Range("C11").Value = 15
$ws.Range("C12").Value = 15
$ws.Range("C13").Value = 15
$ws.Range("C14").Value = 15
$ws.Range("C15").Value = 15
$ws.Range("C16").Value = 15
$ws.Range("C17").Value = 13
$ws.Range("C18").Value = 13
$ws.Range("C19").Value = 13
$ws.Range("C20").Value = 13
$ws.Range("C21").Value = 13
$ws.Range("C22").Value = 13
$ws.Range("C23").Value = 13
$ws.Range("C24").Value = 13
$ws.Range("C25").Value = 14
$ws.Range("C26").Value = 14
$ws.Range("C27").Value = 14
$ws.Range("C28").Value = 12
$ws.Range("C29").Value = 25
$ws.Range("C30").Value = 25
$ws.Range("C31").Value = 25
$ws.Range("C32").Value = 25

# --- Sheet: strategy_id-6005 ---
$ws = $wb.Worksheets.Item("strategy_id-6005")
$ws.Range("C2").Value = 15
$ws.Range("C3").Value = 15
$ws.Range("C4").Value = 15
$ws.Range("C5").Value = 13
$ws.Range("C6").Value = 13
$ws.Range("C7").Value = 13
$ws.Range("C8").Value = 15
$ws.Range("C9").Value = 13
$ws.Range("C10").Value = 13
$ws.Range("C11").Value = 15
$ws.Range("C12").Value = 15
$ws.Range("C13").Value = 15
$ws.Range("C14").Value = 15
$ws.Range("C15").Value = 15
$ws.Range("C16").Value = 15
$ws.Range("C17").Value = 13
$ws.Range("C18").Value = 13
$ws.Range("C19").Value = 13
$ws.Range("C20").Value = 13
$ws.Range("C21").Value = 13
$ws.Range("C22").Value = 13
$ws.Range("C23").Value = 13
$ws.Range("C24").Value = 13
$ws.Range("C25").Value = 14
$ws.Range("C26").Value = 14
$ws.Range("C27").Value = 14
$ws.Range("C28").Value = 12
$ws.Range("C29").Value = 25
$ws.Range("C30").Value = 25
$ws.Range("C31").Value = 25
$ws.Range("C32").Value = 25

# --- Sheet: strategy_id-6006 ---
$ws = $wb.Worksheets.Item("strategy_id-6006")
$ws.Range("C2").Value = 15
$ws.Range("C3").Value = 15
$ws.Range("C4").Value = 15
$ws.Range("C5").Value = 13
$ws.Range("C6").Value = 13
$ws.Range("C7").Value = 13
$ws.Range("C8").Value = 15
$ws.Range("C9").Value = 13
$ws.Range("C10").Value = 13
$ws.Range("C11").Value = 15
$ws.Range("C12").Value = 15
$ws.Range("C13").Value = 15
$ws.Range("C14").Value = 15
$ws.Range("C15").Value = 15
$ws.Range("C16").Value = 15
$ws.Range("C17").Value = 13
$ws.Range("C18").Value = 13
$ws.Range("C19").Value = 13
$ws.Range("C20").Value = 13
$ws.Range("C21").Value = 13
$ws.Range("C22").Value = 13
$ws.Range("C23").Value = 13
$ws.Range("C24").Value = 13
$ws.Range("C25").Value = 14
$ws.Range("C26").Value = 14
$ws.Range("C27").Value = 14
$ws.Range("C28").Value = 12
$ws.Range("C29").Value = 25
$ws.Range("C30").Value = 25
$ws.Range("C31").Value = 25
$ws.Range("C32").Value = 25

# --- Sheet: strategy_id-6007 ---
$ws = $wb.Worksheets.Item("strategy_id-6007")
$ws.Range("C2").Value = 15
$ws.Range("C3").Value = 15
$ws.Range("C4").Value = 15
$ws.Range("C5").Value = 13
$ws.Range("C6").Value = 13
$ws.Range("C7").Value = 13
$ws.Range("C8").Value = 15
$ws.Range("C9").Value = 13
$ws.Range("C10").Value = 13
$ws.Range("C11").Value = 15
$ws.Range("C12").Value = 15
$ws.Range("C13").Value = 15
$ws.Range("C14").Value = 15
$ws.Range("C15").Value = 15
$ws.Range("C16").Value = 15
$ws.Range("C17").Value = 13
$ws.Range("C18").Value = 13
$ws.Range("C19").Value = 13
$ws.Range("C20").Value = 13
$ws.Range("C21").Value = 13
$ws.Range("C22").Value = 13
$ws.Range("C23").Value = 13
$ws.Range("C24").Value = 13
$ws.Range("C25").Value = 14
$ws.Range("C26").Value = 14
$ws.Range("C27").Value = 14
$ws.Range("C28").Value = 12
$ws.Range("C29").Value = 25
$ws.Range("C30").Value = 25
$ws.Range("C31").Value = 25
$ws.Range("C32").Value = 25

# --- Sheet: strategy_id-6011 ---
$ws = $wb.Worksheets.Item("strategy_id-6011")
$ws.Range("C2").Value = 15
$ws.Range("C3").Value = 15
$ws.Range("C4").Value = 15
$ws.Range("C5").Value = 13
$ws.Range("C6").Value = 13
$ws.Range("C7").Value = 13
$ws.Range("C8").Value = 15
$ws.Range("C9").Value = 13
$ws.Range("C10").Value = 13
$ws.Range("C11").Value = 15
$ws.Range("C12").Value = 15
$ws.Range("C13").Value = 15
$ws.Range("C14").Value = 15
$ws.Range("C15").Value = 15
$ws.Range("C16").Value = 15
$ws.Range("C17").Value = 13
$ws.Range("C18").Value = 13
$ws.Range("C19").Value = 13
$ws.Range("C20").Value = 13
$ws.Range("C21").Value = 13
$ws.Range("C22").Value = 13
$ws.Range("C23").Value = 13
$ws.Range("C24").Value = 13
$ws.Range("C25").Value = 14
$ws.Range("C26").Value = 14
$ws.Range("C27").Value = 14
$ws.Range("C28").Value = 12
$ws.Range("C29").Value = 25
$ws.Range("C30").Value = 25
$ws.Range("C31").Value = 25
$ws.Range("C32").Value = 25

# --- Sheet: strategy_id-6015 ---
$ws = $wb.Worksheets.Item("strategy_id-6015")
$ws.Range("C2").Value = 15
$ws.Range("C3").Value = 15
$ws.Range("C4").Value = 15
$ws.Range("C5").Value = 13
$ws.Range("C6").Value = 13
$ws.Range("C7").Value = 13
$ws.Range("C8").Value = 15
$ws.Range("C9").Value = 13
$ws.Range("C10").Value = 13
$ws.Range("C11").Value = 15
$ws.Range("C12").Value = 15
$ws.Range("C13").Value = 15
$ws.Range("C14").Value = 15
$ws.Range("C15").Value = 15
$ws.Range("C16").Value = 15
$ws.Range("C17").Value = 13
$ws.Range("C18").Value = 13
$ws.Range("C19").Value = 13
$ws.Range("C20").Value = 13
$ws.Range("C21").Value = 13
$ws.Range("C22").Value = 13
$ws.Range("C23").Value = 13
$ws.Range("C24").Value = 13
$ws.Range("C25").Value = 14
$ws.Range("C26").Value = 14
$ws.Range("C27").Value = 14
$ws.Range("C28").Value = 12
$ws.Range("C29").Value = 25
$ws.Range("C30").Value = 25
$ws.Range("C31").Value = 25
$ws.Range("C32").Value = 25

# --- Sheet: strategy_id-6016 ---
$ws = $wb.Worksheets.Item("strategy_id-6016")
$ws.Range("C2").Value = 15
$ws.Range("C3").Value = 15
$ws.Range("C4").Value = 15
$ws.Range("C5").Value = 13
$ws.Range("C6").Value = 13
$ws.Range("C7").Value = 13
$ws.Range("C8").Value = 15
$ws.Range("C9").Value = 13
$ws.Range("C10").Value = 13
$ws.Range("C11").Value = 15
$ws.Range("C12").Value = 15
$ws.Range("C13").Value = 15
$ws.Range("C14").Value = 15
$ws.Range("C15").Value = 15
$ws.Range("C16").Value = 15
$ws.Range("C17").Value = 13
$ws.Range("C18").Value = 13
$ws.Range("C19").Value = 13
$ws.Range("C20").Value = 13
$ws.Range("C21").Value = 13
$ws.Range("C22").Value = 13
$ws.Range("C23").Value = 13
$ws.Range("C24").Value = 13
$ws.Range("C25").Value = 14
$ws.Range("C26").Value = 14
$ws.Range("C27").Value = 14
$ws.Range("C28").Value = 12
$ws.Range("C29").Value = 25
$ws.Range("C30").Value = 25
$ws.Range("C31").Value = 25
$ws.Range("C32").Value = 25

